$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) and Volume(1h) (column E) figures for the crypto
# list rows 2-51 to reflect the refreshed values from this run's data pull.
# Cells whose new price text is purely numeric (e.g. "214.51") are forced to
# Text format while writing so Excel keeps them as strings like the source
# data, then restored to the Normal style so formatting is unaffected.
$ws.Range("D2").Value = "25.961.01"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "1.637.69"
$ws.Range("E3").Value = "  -0.17%  "
$ws.Range("E4").Value = "  -0.35%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5087"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.81%  "
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2562"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06340"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.65%  "
$ws.Range("E10").Value = "  +0.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07767"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.35%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.273"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.04%  "
$ws.Range("D13").Value = "1.638.78"
$ws.Range("E13").Value = "  -0.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5416"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "64.01"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.95%  "
$ws.Range("D16").Value = "0.0₅7680"
$ws.Range("E16").Value = "  -2.35%  "
$ws.Range("D17").Value = "25.975.76"
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("E18").Value = "  -0.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "198.94"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.405"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.56%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.882"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.027"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.64%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.869"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.48"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.97%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1195"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.73%  "
$ws.Range("E27").Value = "  -0.34%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.61"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.73%  "
$ws.Range("E29").Value = "  -0.50%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.04889"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.86%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.250"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.163"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.80%  "
$ws.Range("E33").Value = "  -0.41%  "
$ws.Range("E34").Value = "  +0.20%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9040"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.39%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.584"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.81%  "
$ws.Range("D37").Value = "1.139.98"
$ws.Range("E37").Value = "  +0.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5445"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01563"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.36%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.002"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.526"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8088"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.12%  "
$ws.Range("D43").Value = "0.0₈124"
$ws.Range("E43").Value = "  +0.69%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.08"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("E45").Value = "  -4.43%  "
$ws.Range("D46").Value = "1.777.39"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4528"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.04%  "
$ws.Range("E48").Value = "  +0.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.83"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.60%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05118"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.002"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.35%  "
